$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 only contained a stray hyperlink cell (B17, "https://demoqa.com/interaction")
# with no other data -- remove it entirely, shifting TS016..TS019 (old rows 18-21) up
# to rows 17-20.
$ws.Rows("17:17").Delete()

# Deleting the row does not renumber the worksheet-level <hyperlinks> ranges, so rebuild
# them all from scratch (Hyperlinks.Delete() clears every hyperlink on the sheet, not
# just the scoped range, so the full set has to be recreated in order to keep the
# r:id numbering/targets stable).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://demoqa.com/elements")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://demoqa.com/elements")
$ws.Hyperlinks.Add($ws.Range("B4:B8"), "https://demoqa.com/elements", "", "", "https://demoqa.com/elements")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://demoqa.com/forms")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://demoqa.com/alertsWindows")
$ws.Hyperlinks.Add($ws.Range("B11:B12"), "https://demoqa.com/alertsWindows", "", "", "https://demoqa.com/alertsWindows")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://demoqa.com/widgets")
$ws.Hyperlinks.Add($ws.Range("B14:B16"), "https://demoqa.com/widgets", "", "", "https://demoqa.com/widgets")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://demoqa.com/bookstore")
$ws.Hyperlinks.Add($ws.Range("B18:B20"), "https://demoqa.com/bookstore", "", "", "https://demoqa.com/bookstore")

# Hyperlinks.Add stamps the builtin "Hyperlink" style on every cell it touches; restore
# the original column-B look (font/alignment that was already shifted up correctly by
# the row delete) by repainting the formats from an untouched donor cell.
$ws.Range("B16").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's new viewport/selection.
$ws.Range("F24").Select()
